# Scheduled-runner market-data refresh for the Excalibur_Profits workbook.
# For each affected Leve row, re-writes the computed market columns
# (currentAveragePrice.. / LevePriceNQ.. / LeveProfitNQ.. / LeveProfitHQ..,
# columns H:N) with freshly pulled values. A few rows also gain or lose a
# LeveProfitNQ/LeveProfitHQ cell entirely, since those columns are only
# populated when the corresponding craft type has a meaningful profit figure.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 465.8889
$ws.Range("I96").Value = 492.2857
$ws.Range("J96").Value = 373.5
$ws.Range("K96").Value = 1476.8571
$ws.Range("L96").Value = 1120.5
$ws.Range("M96").Value = -103.8571000000002
$ws.Range("N96").Value = -3866.5

$ws.Range("H103").Value = 1163.3334
$ws.Range("I103").Value = 1145
$ws.Range("K103").Value = 3435
$ws.Range("M103").Value = -2849

$ws.Range("H113").Value = 4803.15
$ws.Range("I113").Value = 6128.077
$ws.Range("J113").Value = 2342.5715
$ws.Range("K113").Value = 6128.077
$ws.Range("L113").Value = 2342.5715
$ws.Range("M113").Value = -2874.077
$ws.Range("N113").Value = -8850.5715

$ws.Range("H116").Value = 17093.885
$ws.Range("I116").Value = 16536.291
$ws.Range("J116").Value = 18310.455
$ws.Range("K116").Value = 16536.291
$ws.Range("L116").Value = 18310.455
$ws.Range("M116").Value = -13094.291
$ws.Range("N116").Value = -25194.455

$ws.Range("H127").Value = 2796.1
$ws.Range("I127").Value = 2796.1
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 8388.299999999999
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -3428.299999999999
$ws.Range("N127").ClearContents()

$ws.Range("H129").Value = 2269.75
$ws.Range("I129").Value = 1723.8
$ws.Range("J129").Value = 4999.5
$ws.Range("K129").Value = 5171.4
$ws.Range("L129").Value = 14998.5
$ws.Range("M129").Value = -171.3999999999996
$ws.Range("N129").Value = -24998.5

$ws.Range("H132").Value = 73874.03
$ws.Range("I132").Value = 82563.83
$ws.Range("K132").Value = 247691.49
$ws.Range("M132").Value = -245161.49

$ws.Range("H135").Value = 1201.1177
$ws.Range("I135").Value = 1295.5927
$ws.Range("K135").Value = 11660.3343
$ws.Range("M135").Value = -9125.334299999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4533.6665
$ws.Range("I2").Value = 5621.5713
$ws.Range("J2").Value = 3010.6
$ws.Range("K2").Value = 5621.5713
$ws.Range("L2").Value = 3010.6
$ws.Range("M2").Value = -5508.5713
$ws.Range("N2").Value = -3236.6

$ws.Range("H112").Value = 49864.5
$ws.Range("J112").Value = 49796.75
$ws.Range("L112").Value = 49796.75
$ws.Range("N112").Value = -52750.75

$ws.Range("H116").Value = 4533.6665
$ws.Range("I116").Value = 5621.5713
$ws.Range("J116").Value = 3010.6
$ws.Range("K116").Value = 5621.5713
$ws.Range("L116").Value = 3010.6
$ws.Range("M116").Value = -3327.5713
$ws.Range("N116").Value = -7598.6

$ws.Range("H132").Value = 246809.73
$ws.Range("I132").Value = 418763.22
$ws.Range("J132").Value = 6074.85
$ws.Range("K132").Value = 1256289.66
$ws.Range("L132").Value = 18224.55
$ws.Range("M132").Value = -1253759.66
$ws.Range("N132").Value = -23284.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4533.6665
$ws.Range("I3").Value = 5621.5713
$ws.Range("J3").Value = 3010.6
$ws.Range("K3").Value = 5621.5713
$ws.Range("L3").Value = 3010.6
$ws.Range("M3").Value = -5507.5713
$ws.Range("N3").Value = -3238.6

$ws.Range("H86").Value = 2500
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1377
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 2500
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6884
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 743.8387
$ws.Range("I94").Value = 613.3333
$ws.Range("J94").Value = 1624.75
$ws.Range("K94").Value = 613.3333
$ws.Range("L94").Value = 1624.75
$ws.Range("M94").Value = -162.3333
$ws.Range("N94").Value = -2526.75

$ws.Range("H107").Value = 1935.48
$ws.Range("I107").Value = 2253.6667
$ws.Range("J107").Value = 1117.2858
$ws.Range("K107").Value = 2253.6667
$ws.Range("L107").Value = 1117.2858
$ws.Range("M107").Value = -333.6667000000002
$ws.Range("N107").Value = -4957.2858

$ws.Range("H108").Value = 98684
$ws.Range("J108").Value = 98684
$ws.Range("L108").Value = 98684
$ws.Range("N108").Value = -106364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 970.5454999999999
$ws.Range("I16").Value = 990.75
$ws.Range("J16").Value = 916.6667
$ws.Range("K16").Value = 990.75
$ws.Range("L16").Value = 916.6667
$ws.Range("M16").Value = -703.75
$ws.Range("N16").Value = -1490.6667

$ws.Range("H22").Value = 597.0769
$ws.Range("I22").Value = 507
$ws.Range("J22").Value = 897.3333
$ws.Range("K22").Value = 507
$ws.Range("L22").Value = 897.3333
$ws.Range("M22").Value = -157
$ws.Range("N22").Value = -1597.3333

$ws.Range("H31").Value = 7747.442
$ws.Range("I31").Value = 1940.3334
$ws.Range("J31").Value = 10858.393
$ws.Range("K31").Value = 1940.3334
$ws.Range("L31").Value = 10858.393
$ws.Range("M31").Value = -1645.3334
$ws.Range("N31").Value = -11448.393

$ws.Range("H34").Value = 7747.442
$ws.Range("I34").Value = 1940.3334
$ws.Range("J34").Value = 10858.393
$ws.Range("K34").Value = 1940.3334
$ws.Range("L34").Value = 10858.393
$ws.Range("M34").Value = -1738.3334
$ws.Range("N34").Value = -11262.393

$ws.Range("H47").Value = 11000
$ws.Range("I47").Value = 11000
$ws.Range("K47").Value = 11000
$ws.Range("M47").Value = -10434

$ws.Range("H113").Value = 970.5454999999999
$ws.Range("I113").Value = 990.75
$ws.Range("J113").Value = 916.6667
$ws.Range("K113").Value = 990.75
$ws.Range("L113").Value = 916.6667
$ws.Range("M113").Value = 1179.25
$ws.Range("N113").Value = -5256.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 609.8333
$ws.Range("I14").Value = 609.8333
$ws.Range("K14").Value = 1829.4999
$ws.Range("M14").Value = -1656.4999

$ws.Range("H41").Value = 132.66667
$ws.Range("I41").Value = 99
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 297
$ws.Range("L41").Value = 600
$ws.Range("M41").Value = 41
$ws.Range("N41").Value = -1276

$ws.Range("H131").Value = 2633.5557
$ws.Range("I131").Value = 10349.5
$ws.Range("K131").Value = 31048.5
$ws.Range("M131").Value = -26008.5

$ws.Range("H132").Value = 2514.2778
$ws.Range("I132").Value = 2199.4167
$ws.Range("J132").Value = 2671.7083
$ws.Range("K132").Value = 19794.7503
$ws.Range("L132").Value = 24045.3747
$ws.Range("M132").Value = -17264.7503
$ws.Range("N132").Value = -29105.3747

$ws.Range("H137").Value = 4823.4736
$ws.Range("J137").Value = 9195.286
$ws.Range("L137").Value = 27585.858
$ws.Range("N137").Value = -37785.858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3824.75
$ws.Range("J102").Value = 5315.25
$ws.Range("L102").Value = 5315.25
$ws.Range("N102").Value = -8559.25

$ws.Range("H107").Value = 478.8125
$ws.Range("J107").Value = 748.8333
$ws.Range("L107").Value = 748.8333
$ws.Range("N107").Value = -4588.8333

$ws.Range("H111").Value = 59999
$ws.Range("J111").Value = 59999
$ws.Range("L111").Value = 59999
$ws.Range("N111").Value = -66133

$ws.Range("H132").Value = 258196.45
$ws.Range("I132").Value = 378232.06
$ws.Range("J132").Value = 2120.4666
$ws.Range("K132").Value = 1134696.18
$ws.Range("L132").Value = 6361.399800000001
$ws.Range("M132").Value = -1132166.18
$ws.Range("N132").Value = -11421.3998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 23577.092
$ws.Range("I61").Value = 25745.3
$ws.Range("J61").Value = 1895
$ws.Range("K61").Value = 25745.3
$ws.Range("L61").Value = 1895
$ws.Range("M61").Value = -25543.3
$ws.Range("N61").Value = -2299

$ws.Range("H110").Value = 59790.777
$ws.Range("J110").Value = 59790.777
$ws.Range("L110").Value = 59790.777
$ws.Range("N110").Value = -67970.777

$ws.Range("H113").Value = 23577.092
$ws.Range("I113").Value = 25745.3
$ws.Range("J113").Value = 1895
$ws.Range("K113").Value = 25745.3
$ws.Range("L113").Value = 1895
$ws.Range("M113").Value = -23575.3
$ws.Range("N113").Value = -6235
